$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This edit rotates the species-observation data among rows 2, 3 and 6:
#   new row 2 <- old row 3 data
#   new row 3 <- old row 6 data (plus the extra "blomning/Blommande" fields)
#   new row 6 <- old row 2 data (clearing the extra fields row 6 used to carry)

# --- Row 2 gets the old Row 3 values ---
$ws.Range("A2").Value = 111525999
$ws.Range("B2").Value = 86223
$ws.Range("D2").Value = "NT"
$ws.Range("E2").Value = 4412
$ws.Range("F2").Value = "Äggvaxskivling"
$ws.Range("G2").Value = "Hygrophorus karstenii"
$ws.Range("H2").Value = "Sacc. & Cub."
$ws.Range("Q2").Value = 538398.3112996884
$ws.Range("R2").Value = 7024277.647416403

# --- Row 3 gets the old Row 6 values (including the bloom-stage fields) ---
$ws.Range("A3").Value = 111526007
$ws.Range("B3").Value = 96348
$ws.Range("D3").Value = "VU"
$ws.Range("E3").Value = 220787
$ws.Range("F3").Value = "Knärot"
$ws.Range("G3").Value = "Goodyera repens"
$ws.Range("H3").Value = "(L.) R. Br."
$ws.Range("J3").Value = ""
$ws.Range("K3").Value = "blomning"
$ws.Range("L3").Value = ""
$ws.Range("N3").Value = ""
$ws.Range("Q3").Value = 538522.0815204142
$ws.Range("R3").Value = 7024306.075093818
$ws.Range("AC3").Value = "Blommande"
$ws.Range("AF3").Value = ""

# --- Row 6 gets the old Row 2 values (and loses the bloom-stage fields) ---
$ws.Range("A6").Value = 111525965
$ws.Range("B6").Value = 90332
$ws.Range("D6").Value = "LC"
$ws.Range("E6").Value = 4769
$ws.Range("F6").Value = "Svavelriska"
$ws.Range("G6").Value = "Lactarius scrobiculatus"
$ws.Range("H6").Value = "(Scop.:Fr.) Fr."
$ws.Range("K6").ClearContents()
$ws.Range("Q6").Value = 538536.3052666293
$ws.Range("R6").Value = 7024282.445275509
$ws.Range("AC6").ClearContents()
